$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue ($ws.Range("D2")) "308.12"
Set-TextValue ($ws.Range("E2")) "1.83%"

# Row 3
Set-TextValue ($ws.Range("D3")) "37.82"
Set-TextValue ($ws.Range("E3")) "5.80%"

# Row 4
Set-TextValue ($ws.Range("D4")) "5.066"
Set-TextValue ($ws.Range("E4")) "0.65%"

# Row 5
Set-TextValue ($ws.Range("D5")) "0.08146"
Set-TextValue ($ws.Range("E5")) "3.20%"

# Row 6
Set-TextValue ($ws.Range("D6")) "1.975"
Set-TextValue ($ws.Range("E6")) "6.89%"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue ($ws.Range("D7")) "4.175"
Set-TextValue ($ws.Range("E7")) "1.55%"

# Row 8
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue ($ws.Range("D8")) "7.899"
Set-TextValue ($ws.Range("E8")) "1.57%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue ($ws.Range("D9")) "0.9269"
Set-TextValue ($ws.Range("E9")) "0.17%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue ($ws.Range("D10")) "0.1392"
Set-TextValue ($ws.Range("E10")) "2.67%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue ($ws.Range("D11")) "0.1944"
Set-TextValue ($ws.Range("E11")) "2.23%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue ($ws.Range("D12")) "0.09247"
Set-TextValue ($ws.Range("E12")) "1.47%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue ($ws.Range("D13")) "0.03509"
Set-TextValue ($ws.Range("E13")) "1.02%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue ($ws.Range("D14")) "0.09885"
Set-TextValue ($ws.Range("E14")) "0.50%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue ($ws.Range("D15")) "0.001407"
Set-TextValue ($ws.Range("E15")) "0.27%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue ($ws.Range("D16")) "0.006238"
Set-TextValue ($ws.Range("E16")) "0.59%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue ($ws.Range("D17")) "3.941"
Set-TextValue ($ws.Range("E17")) "5.87%"

# Row 18
Set-TextValue ($ws.Range("D18")) "3.409"
Set-TextValue ($ws.Range("E18")) "2.07%"

# Row 19
Set-TextValue ($ws.Range("D19")) "0.3455"
Set-TextValue ($ws.Range("E19")) "0.39%"

# Row 20
Set-TextValue ($ws.Range("D20")) "0.1285"
Set-TextValue ($ws.Range("E20")) "-4.34%"

# Row 21
Set-TextValue ($ws.Range("D21")) "4.803"
Set-TextValue ($ws.Range("E21")) "-6.93%"

# Row 22
Set-TextValue ($ws.Range("D22")) "0.2619"
Set-TextValue ($ws.Range("E22")) "18.86%"

# Row 23
Set-TextValue ($ws.Range("D23")) "0.04435"
Set-TextValue ($ws.Range("E23")) "0.49%"

# Row 24
Set-TextValue ($ws.Range("D24")) "0.001245"
Set-TextValue ($ws.Range("E24")) "0.82%"

# Row 25
Set-TextValue ($ws.Range("E25")) "-9.79%"

# Row 27
Set-TextValue ($ws.Range("D27")) "0.0001300"
Set-TextValue ($ws.Range("E27")) "-0.17%"

# Row 39
Set-TextValue ($ws.Range("D39")) "0.02115"
Set-TextValue ($ws.Range("E39")) "9.37%"

# Row 40
Set-TextValue ($ws.Range("D40")) "0.05141"
Set-TextValue ($ws.Range("E40")) "1.06%"

# Row 41
Set-TextValue ($ws.Range("D41")) "0.007476"
Set-TextValue ($ws.Range("E41")) "-0.98%"

# Row 42
Set-TextValue ($ws.Range("D42")) "0.01013"

# Row 43
Set-TextValue ($ws.Range("E43")) "1.42%"

# Row 44
Set-TextValue ($ws.Range("D44")) "0.002130"
Set-TextValue ($ws.Range("E44")) "-1.11%"

# Row 45
Set-TextValue ($ws.Range("E45")) "-0.46%"

# Row 46
Set-TextValue ($ws.Range("D46")) "0.00006389"
Set-TextValue ($ws.Range("E46")) "4.00%"

# Row 47
Set-TextValue ($ws.Range("E47")) "-0.15%"

# Row 48
Set-TextValue ($ws.Range("E48")) "1.94%"

# Row 49
Set-TextValue ($ws.Range("D49")) "0.001603"
Set-TextValue ($ws.Range("E49")) "-3.54%"

# Row 50
Set-TextValue ($ws.Range("D50")) "0.00002100"
Set-TextValue ($ws.Range("E50")) "-0.15%"

# Row 51
Set-TextValue ($ws.Range("D51")) "0.0002000"
Set-TextValue ($ws.Range("E51")) "-0.15%"
